$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 309
$ws.Range("I41").Value = 96.666664
$ws.Range("J41").Value = 436.4
$ws.Range("K41").Value = 96.666664
$ws.Range("L41").Value = 436.4
$ws.Range("M41").Value = 343.333336
$ws.Range("N41").Value = -1316.4
$ws.Range("H76").Value = 3369.889
$ws.Range("I76").Value = 3265.8
$ws.Range("K76").Value = 3265.8
$ws.Range("M76").Value = -2950.8
$ws.Range("H79").Value = 3369.889
$ws.Range("I79").Value = 3265.8
$ws.Range("K79").Value = 3265.8
$ws.Range("M79").Value = -2173.8
$ws.Range("H86").Value = 13403.125
$ws.Range("I86").Value = 1040.6
$ws.Range("J86").Value = 34007.332
$ws.Range("K86").Value = 1040.6
$ws.Range("L86").Value = 34007.332
$ws.Range("M86").Value = 82.40000000000009
$ws.Range("N86").Value = -36253.332
$ws.Range("H89").Value = 13403.125
$ws.Range("I89").Value = 1040.6
$ws.Range("J89").Value = 34007.332
$ws.Range("K89").Value = 5203
$ws.Range("L89").Value = 170036.66
$ws.Range("M89").Value = 413
$ws.Range("N89").Value = -181268.66
$ws.Range("H98").Value = 599.5333000000001
$ws.Range("I98").Value = 346.3846
$ws.Range("J98").Value = 2245
$ws.Range("K98").Value = 346.3846
$ws.Range("L98").Value = 2245
$ws.Range("M98").Value = 1151.6154
$ws.Range("N98").Value = -5241
$ws.Range("H122").Value = 599.5333000000001
$ws.Range("I122").Value = 346.3846
$ws.Range("J122").Value = 2245
$ws.Range("K122").Value = 1039.1538
$ws.Range("L122").Value = 6735
$ws.Range("M122").Value = 1410.8462
$ws.Range("N122").Value = -11635
$ws.Range("H129").Value = 205941.14
$ws.Range("J129").Value = 210221.17
$ws.Range("L129").Value = 630663.51
$ws.Range("N129").Value = -640663.51
$ws.Range("H137").Value = 41516.84
$ws.Range("I137").Value = 1034.4445
$ws.Range("K137").Value = 3103.3335
$ws.Range("M137").Value = -553.3335000000002
$ws.Range("H141").Value = 1932.7826
$ws.Range("I141").Value = 1365.3125
$ws.Range("J141").Value = 3229.8572
$ws.Range("K141").Value = 4095.9375
$ws.Range("L141").Value = 9689.571599999999
$ws.Range("M141").Value = 1084.0625
$ws.Range("N141").Value = -20049.5716
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2461.3784
$ws.Range("I45").Value = 2429.625
$ws.Range("K45").Value = 2429.625
$ws.Range("M45").Value = -2052.625
$ws.Range("H61").Value = 1322.4348
$ws.Range("I61").Value = 1347.0476
$ws.Range("J61").Value = 1064
$ws.Range("K61").Value = 1347.0476
$ws.Range("L61").Value = 1064
$ws.Range("M61").Value = -1135.0476
$ws.Range("N61").Value = -1488
$ws.Range("H74").Value = 38463904
$ws.Range("J74").Value = 928.5
$ws.Range("L74").Value = 928.5
$ws.Range("N74").Value = -2676.5
$ws.Range("H77").Value = 38463904
$ws.Range("J77").Value = 928.5
$ws.Range("L77").Value = 4642.5
$ws.Range("N77").Value = -13378.5
$ws.Range("H136").Value = 1322.4348
$ws.Range("I136").Value = 1347.0476
$ws.Range("J136").Value = 1064
$ws.Range("K136").Value = 4041.142800000001
$ws.Range("L136").Value = 3192
$ws.Range("M136").Value = -1491.142800000001
$ws.Range("N136").Value = -8292
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 500317.16
$ws.Range("I22").Value = 666947.2
$ws.Range("J22").Value = 427
$ws.Range("K22").Value = 666947.2
$ws.Range("L22").Value = 427
$ws.Range("M22").Value = -666774.2
$ws.Range("N22").Value = -773
$ws.Range("H86").Value = 1635.3549
$ws.Range("I86").Value = 1409
$ws.Range("J86").Value = 2110.7
$ws.Range("K86").Value = 1409
$ws.Range("L86").Value = 2110.7
$ws.Range("M86").Value = -286
$ws.Range("N86").Value = -4356.7
$ws.Range("H89").Value = 1635.3549
$ws.Range("I89").Value = 1409
$ws.Range("J89").Value = 2110.7
$ws.Range("K89").Value = 7045
$ws.Range("L89").Value = 10553.5
$ws.Range("M89").Value = -1429
$ws.Range("N89").Value = -21785.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8619.49
$ws.Range("I31").Value = 9335.553
$ws.Range("J31").Value = 6145.8184
$ws.Range("K31").Value = 9335.553
$ws.Range("L31").Value = 6145.8184
$ws.Range("M31").Value = -9040.553
$ws.Range("N31").Value = -6735.8184
$ws.Range("H34").Value = 8619.49
$ws.Range("I34").Value = 9335.553
$ws.Range("J34").Value = 6145.8184
$ws.Range("K34").Value = 9335.553
$ws.Range("L34").Value = 6145.8184
$ws.Range("M34").Value = -9133.553
$ws.Range("N34").Value = -6549.8184
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1547.2609
$ws.Range("J5").Value = 1415.1666
$ws.Range("L5").Value = 4245.4998
$ws.Range("N5").Value = -4469.4998
$ws.Range("H6").Value = 73.35714
$ws.Range("I6").Value = 52.25
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 156.75
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -43.75
$ws.Range("N6").Value = -826
$ws.Range("H135").Value = 1547.2609
$ws.Range("J135").Value = 1415.1666
$ws.Range("L135").Value = 12736.4994
$ws.Range("N135").Value = -17806.4994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31253094
$ws.Range("I102").Value = 31253094
$ws.Range("K102").Value = 31253094
$ws.Range("M102").Value = -31251472
$ws.Range("H126").Value = 4012.487
$ws.Range("I126").Value = 3219.3704
$ws.Range("K126").Value = 9658.111199999999
$ws.Range("M126").Value = -7188.111199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5759.353
$ws.Range("I7").Value = 3210.7
$ws.Range("J7").Value = 9400.286
$ws.Range("K7").Value = 3210.7
$ws.Range("L7").Value = 9400.286
$ws.Range("M7").Value = -3098.7
$ws.Range("N7").Value = -9624.286
$ws.Range("H40").Value = 5400.625
$ws.Range("I40").Value = 6100
$ws.Range("J40").Value = 4701.25
$ws.Range("K40").Value = 6100
$ws.Range("L40").Value = 4701.25
$ws.Range("M40").Value = -5964
$ws.Range("N40").Value = -4973.25
$ws.Range("H82").Value = 2554.7273
$ws.Range("I82").Value = 2266.889
$ws.Range("J82").Value = 3850
$ws.Range("K82").Value = 2266.889
$ws.Range("L82").Value = 3850
$ws.Range("M82").Value = -1905.889
$ws.Range("N82").Value = -4572
$ws.Range("H85").Value = 2554.7273
$ws.Range("I85").Value = 2266.889
$ws.Range("J85").Value = 3850
$ws.Range("K85").Value = 2266.889
$ws.Range("L85").Value = 3850
$ws.Range("M85").Value = -1018.889
$ws.Range("N85").Value = -6346
$ws.Range("H126").Value = 5759.353
$ws.Range("I126").Value = 3210.7
$ws.Range("J126").Value = 9400.286
$ws.Range("K126").Value = 9632.099999999999
$ws.Range("L126").Value = 28200.858
$ws.Range("M126").Value = -7162.099999999999
$ws.Range("N126").Value = -33140.858
$ws.Range("H135").Value = 50429
$ws.Range("J135").Value = 50429
$ws.Range("L135").Value = 50429
$ws.Range("N135").Value = -60569
$ws.Range("H136").Value = 30995.53
$ws.Range("I136").Value = 39908.383
$ws.Range("J136").Value = 2028.75
$ws.Range("K136").Value = 119725.149
$ws.Range("L136").Value = 6086.25
$ws.Range("M136").Value = -117175.149
$ws.Range("N136").Value = -11186.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 74357.14
$ws.Range("I14").Value = 100520
$ws.Range("J14").Value = 63892
$ws.Range("K14").Value = 100520
$ws.Range("L14").Value = 63892
$ws.Range("M14").Value = -100352
$ws.Range("N14").Value = -64228
$ws.Range("H96").Value = 1503.2
$ws.Range("I96").Value = 1172
$ws.Range("K96").Value = 1172
$ws.Range("M96").Value = 201
$ws.Range("H122").Value = 1275.1111
$ws.Range("I122").Value = 1360.2
$ws.Range("J122").Value = 849.6667
$ws.Range("K122").Value = 4080.6
$ws.Range("L122").Value = 2549.0001
$ws.Range("M122").Value = -1630.6
$ws.Range("N122").Value = -7449.0001
$ws.Range("H126").Value = 1251.0834
$ws.Range("I126").Value = 1125.75
$ws.Range("J126").Value = 1501.75
$ws.Range("K126").Value = 3377.25
$ws.Range("L126").Value = 4505.25
$ws.Range("M126").Value = -907.25
$ws.Range("N126").Value = -9445.25
